$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "line" value for the setDeviceInitializationContext row.
$ws.Range("D2").Value = 1376

# Bring D3:D4's direct formatting in line with D2's (copy format only,
# leaving their numeric values untouched).
$ws.Range("D2").Copy()
$ws.Range("D3:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to C12, matching the saved cursor position.
$ws.Range("C12").Select()
